# Update "想去人数" (want-to-go count) values in column F for both the
# "展览" (Exhibition) sheet and the "全部类型" (All Types) sheet, which
# contain duplicate data.

$wb = $excel.ActiveWorkbook

# Row => new value for column F
$updates = @{
    7  = 107
    8  = 79
    11 = 19
    14 = 291
    16 = 352
    20 = 43
    22 = 872
    24 = 291
    30 = 83
    31 = 206
    33 = 263
    34 = 1601
    40 = 3509
    42 = 183
    43 = 891
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
